# Weekly update: insert a new weekly price-report row for "Poroto verde"
# (Agrícola del Norte S.A. de Arica) at row 8, pushing the existing
# historical rows (old 8..45) down to (9..46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 8 - this shifts rows 8:45 down to 9:46
# and also extends the used range / dimension to row 46 automatically.
$ws.Rows("8:8").Insert()

# Populate the newly-inserted row 8 with this week's data point.
$ws.Cells.Item(8, 1).Value  = 1
$ws.Cells.Item(8, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value  = 44547
$ws.Cells.Item(8, 5).Value  = 15
$ws.Cells.Item(8, 6).Value  = 100112031
$ws.Cells.Item(8, 7).Value  = "Poroto verde"
$ws.Cells.Item(8, 8).Value  = "Sin especificar"
$ws.Cells.Item(8, 9).Value  = "Primera"
$ws.Cells.Item(8, 10).Value = 1700
$ws.Cells.Item(8, 11).Value = 400
$ws.Cells.Item(8, 12).Value = 500
$ws.Cells.Item(8, 13).Value = 450
$ws.Cells.Item(8, 14).Value = "$/kilo"
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 450
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = "Hortaliza"
